# "Generate Report for Handback"
#
# For both locale sheets (zh-cn, de-de) this adds two new report columns
# per data row: the "Latest Target File" (F) and "Latest Handback File" (G)
# hyperlinked cells, mirroring the existing "Latest Handoff File" (A) /
# "Latest Target File" (D) hyperlinks but pointing at the canonical
# source/handoff-target pair. It also updates the Status text (now that the
# file has been handed back) and records the handback timestamp in column H.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
#    This shared string is referenced from every sheet (Overview, zh-cn,
#    de-de), so a single text replace keeps every occurrence in sync.
# ---------------------------------------------------------------------
foreach ($sheetName in @("Overview", "zh-cn", "de-de")) {
    $sheet = $wb.Worksheets.Item($sheetName)
    $sheet.Cells.Replace("Ready for handoff", "Handed back: in sync with en-US")
}

# ---------------------------------------------------------------------
# 2) Per-locale sheet updates: add F/G "Latest Target File" / "Latest
#    Handback File" hyperlinked cells for rows 2 and 3, and stamp the
#    handback datetime into column H.
# ---------------------------------------------------------------------
$mdFileName = "567f5b35-998d-4faf-aa16-a2c0f7cb7d55.md"
$mdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/513719c527e3fbced7892f685fe750713cbf4aba/e2e/567f5b35-998d-4faf-aa16-a2c0f7cb7d55.md"

$locales = @(
    @{
        Sheet = "zh-cn"
        XlfFileName = "567f5b35-998d-4faf-aa16-a2c0f7cb7d55.d0dc26581f56af7f0997c05bfe685f1717bc0c37.zh-cn.xlf"
        XlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b8c69a350b6334dbdb53566eebfb12e567f60739/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/567f5b35-998d-4faf-aa16-a2c0f7cb7d55.d0dc26581f56af7f0997c05bfe685f1717bc0c37.zh-cn.xlf"
        HandbackDateTime = "2016-03-21 15:04:25"
    },
    @{
        Sheet = "de-de"
        XlfFileName = "567f5b35-998d-4faf-aa16-a2c0f7cb7d55.d0dc26581f56af7f0997c05bfe685f1717bc0c37.de-de.xlf"
        XlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4ec60a443d7ffa5d9757a0811fcaaec9e881152f/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/567f5b35-998d-4faf-aa16-a2c0f7cb7d55.d0dc26581f56af7f0997c05bfe685f1717bc0c37.de-de.xlf"
        HandbackDateTime = "2016-03-21 15:04:31"
    }
)

foreach ($locale in $locales) {
    $ws = $wb.Worksheets.Item($locale.Sheet)

    foreach ($row in @(2, 3)) {
        $fCell = $ws.Range("F" + $row)
        $gCell = $ws.Range("G" + $row)

        $fCell.Value = $mdFileName
        $ws.Hyperlinks.Add($fCell, $mdUrl, "", "", $mdFileName) | Out-Null
        $fCell.Style = "HyperLink"

        $gCell.Value = $locale.XlfFileName
        $ws.Hyperlinks.Add($gCell, $locale.XlfUrl, "", "", $locale.XlfFileName) | Out-Null
        $gCell.Style = "HyperLink"

        $ws.Range("H" + $row).Value = $locale.HandbackDateTime
    }
}
